# Replace the final "조합 분할" paragraph with itself (now un-styled, pPr removed)
# plus the newly authored "046 분산 데이터베이스 설계" and "047 데이터베이스 이중화 / 서버
# 클러스터링" sections that were appended after it (2장 046~047).
$d = $word.ActiveDocument

# Locate the paragraph that currently ends the document ("조합 분할 – 범위 분할로 ...").
$searchRng = $d.Content
$found = $searchRng.Find.Execute("조합 분할 – 범위 분할로 분할한 다음 해시 함수를 적용하여 다시 분할하는 방식", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "could not locate the target paragraph to anchor the edit"
}
$target = $searchRng.Paragraphs(1).Range

# InsertXML replaces the addressed range's contents outright, so handing it the whole
# replacement run of paragraphs (old paragraph, minus its pPr, followed by the 18 new
# paragraphs through the new RPO paragraph + the relocated _GoBack bookmark) reproduces
# the diff in a single call.
$xmlFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">조합 분할 </w:t></w:r><w:r><w:t xml:space="preserve">– </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">범위 분할로 분할한 다음 해시 함수를 적용하여 다시 분할하는 방식 </w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:sz w:val="22"/><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="22"/><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">046 </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>분산 데이터베이스 설계</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">데이터베이스 용량 설계 </w:t></w:r><w:r><w:t xml:space="preserve">– </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>데이터가 저장될 공간을 정의하는 것</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">테이블에 저장될 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>데이터양과</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> 인덱스,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>클러스터 등이 차지하는 공간 등을 예측하여 반영해야 함</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">분산 데이터베이스 </w:t></w:r><w:r><w:t xml:space="preserve">– </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>논리적으로는 하나의 시스템에 속하지만 물리적으로는 네트워크를 통해 연결된 여러 개의 사이트에 분산된 데이터베이스를 말함</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">분산 데이터베이스의 목표 </w:t></w:r><w:r><w:t xml:space="preserve">– </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">위치 투명성 </w:t></w:r><w:r><w:t xml:space="preserve">/ </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">중복 투명성 / 병행 투명성 </w:t></w:r><w:r><w:t xml:space="preserve">/ </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>장애 투명성</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">분산 설계 방법 </w:t></w:r><w:r><w:t>–</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> 테이블 위치 분산 </w:t></w:r><w:r><w:t xml:space="preserve">/ </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">분할 </w:t></w:r><w:r><w:t xml:space="preserve">/ </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>할당</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:sz w:val="22"/><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">047 데이터베이스 이중화 </w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">/ </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">서버 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>클러스터링</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">데이터베이스 이중화 </w:t></w:r><w:r><w:t xml:space="preserve">– </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>시스템 오류로 인한 데이터베이스 서비스 중단이나 물리적 손상 발생 시 복구하기 위해 동일한 데이터베이스를 복제하여 관리하는 것</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>- 하나 이상의 데이터베이스가 항상 같은 상태를 유지하므로 데이터베이스에 문제가 발생하면 복제된 데이터베이스를 이용하여 즉시 문제를 해결할 수 있음</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">데이터베이스 이중화의 분류 </w:t></w:r><w:r><w:t xml:space="preserve">– Eager </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">기법 </w:t></w:r><w:r><w:t xml:space="preserve">/ Lazy </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>기법</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">데이터베이스 이중화 구성 방법 </w:t></w:r><w:r><w:t xml:space="preserve">– </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">활동-대기 방법 </w:t></w:r><w:r><w:t xml:space="preserve">/ </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>활동-활동 방법</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>클러스터링</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">– </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>두 대 이상의 서버를 하나의 서버처럼 운영하는 기술</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>클러스터링</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> 종류 </w:t></w:r><w:r><w:t xml:space="preserve">– </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>고가용성</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>클러스터링</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> /</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">병렬 처리 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>클러스터링</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t>RTO</w:t></w:r><w:r><w:t>(</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>목표 복구 시간)</w:t></w:r><w:r><w:t xml:space="preserve"> – </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>비상상태 또는 업무 중단 시점으로부터 복구되어 가동될 때까지의 소요 시간을 의미</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">RPO - </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>비상상태 또는 업무 중단 시점으로부터</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> 데이터를 복구할 수 있는 기준점을 의미</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$target.InsertXML($xmlFrag)
